# Fruta / hortaliza, semanal
# Insert a new week (44578) of Melon price data above the existing last
# week (44571) entries for "Macroferia Regional de Talca".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows before row 268; this pushes the current rows
# 268:273 (the 44571 week, 6 line items) down to 273:278.
$ws.Rows("268:272").Insert()

# New week's data (date 44578) - 5 line items.
$rows = @(
    @{ Row = 268; Variedad = "Calameño"; Calidad = "Primera"; Volumen = 5000; Precio = 700 },
    @{ Row = 269; Variedad = "Calameño"; Calidad = "Segunda"; Volumen = 3000; Precio = 500 },
    @{ Row = 270; Variedad = "Tuna";     Calidad = "Primera"; Volumen = 4000; Precio = 700 },
    @{ Row = 271; Variedad = "Tuna";     Calidad = "Segunda"; Volumen = 3500; Precio = 500 },
    @{ Row = 272; Variedad = "Tuna";     Calidad = "Tercera"; Volumen = 2000; Precio = 300 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 5
    $ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($row, 3).Value = "Maule"
    $ws.Cells.Item($row, 4).Value = 44578
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = 7
    $ws.Cells.Item($row, 6).Value = 100112027
    $ws.Cells.Item($row, 7).Value = "Melón"
    $ws.Cells.Item($row, 8).Value = $r.Variedad
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.Precio
    $ws.Cells.Item($row, 12).Value = $r.Precio
    $ws.Cells.Item($row, 13).Value = $r.Precio
    $ws.Cells.Item($row, 14).Value = "`$/unidad"
    $ws.Cells.Item($row, 15).Value = "Región del Maule"
    $ws.Cells.Item($row, 16).Value = $r.Precio
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
